# Implement FastAPI backend for AI-Based PowerPoint
# Replaces the generic "Sample Presentation" placeholder deck text with the
# real sample.txt-derived content (titles + bullet content on slides 2 & 3).

$p = $ppt.ActivePresentation

# Helper: replace a paragraph's text while avoiding the engine's
# sentence-boundary run-split quirk that fires when new text ends with a
# period (clearing first, then inserting, keeps it a single run).
function Set-ParaText($para, [string]$text) {
    $para.Text = ""
    $para.InsertAfter($text) | Out-Null
}

# ---------------------------------------------------------------------
# Slide 1: title card
# ---------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$tr1Title = $s1.Shapes.Item(1).TextFrame.TextRange
Set-ParaText $tr1Title.Paragraphs(1, 1) "sample.txt"

# ---------------------------------------------------------------------
# Slide 2: "(part 1)" title + first five bullet paragraphs
# ---------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$tr2Title = $s2.Shapes.Item(1).TextFrame.TextRange
Set-ParaText $tr2Title.Paragraphs(1, 1) "sample.txt(part 1)"

$tr2 = $s2.Shapes.Item(2).TextFrame.TextRange
Set-ParaText $tr2.Paragraphs(2, 1) "*   **Pervasiveness of Technology:** Integrated into all aspects of daily life (communication, work, learning)."
Set-ParaText $tr2.Paragraphs(3, 1) "*   **Rapid Innovation:** Evolution from simple machines to advanced AI transforms human potential and simplifies complex tasks."
Set-ParaText $tr2.Paragraphs(4, 1) "*   **The Need for Responsibility:** Technological growth introduces critical challenges, including privacy concerns, misinformation, and digital addiction."
Set-ParaText $tr2.Paragraphs(5, 1) "*   **Ethical Balance:** Societies must balance convenience with consciousness."
Set-ParaText $tr2.Paragraphs(6, 1) "*   **Guiding Principle:** Technology must be used as an *enabler*, not a *controller*."

# ---------------------------------------------------------------------
# Slide 3: "(part 2)" title + bullets collapse from five points to one
# ---------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$tr3Title = $s3.Shapes.Item(1).TextFrame.TextRange
Set-ParaText $tr3Title.Paragraphs(1, 1) "sample.txt(part 2)"

$tr3 = $s3.Shapes.Item(2).TextFrame.TextRange
# Drop the now-obsolete "Point 6..Point 9" paragraphs first (always removing
# paragraph index 2, i.e. everything except the leading blank paragraph and
# the final "Point 10" one). Deleting the text frame's *current last*
# paragraph leaves a stray empty <a:p/> behind in this engine, so we
# deliberately keep "Point 10" alive until the end and simply retarget its
# text afterwards instead of deleting + re-inserting a fresh paragraph.
for ($i = 0; $i -lt 4; $i++) {
    $tr3.Paragraphs(2, 1).Delete()
}
Set-ParaText $tr3.Paragraphs(2, 1) "*   **Future Focus:** Innovation must be guided by empathy, ethics, and purpose to create a kinder and more sustainable world."
